# Apply crypto price/volume refresh as scraped by GitHub Actions job
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    # Force text storage (avoid Excel auto-converting dotted numbers like
    # "52.163.75" or trimming trailing zeros like "0.950") then restore the
    # original (default) cell style so formatting is left untouched.
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" '52.163.75'
Set-TextValue "E2" '  +1.18%  '

Set-TextValue "D3" '2.893.31'
Set-TextValue "E3" '  +3.58%  '

Set-TextValue "E4" '  +0.06%  '

Set-TextValue "D5" '352.20'
Set-TextValue "E5" '  -0.12%  '

Set-TextValue "D6" '112.04'
Set-TextValue "E6" '  +2.97%  '

Set-TextValue "E7" '  +1.58%  '

Set-TextValue "E8" '  +0.10%  '

Set-TextValue "D9" '0.623'
Set-TextValue "E9" '  -0.37%  '

Set-TextValue "E10" '  +1.01%  '

Set-TextValue "D12" '0.0859'
Set-TextValue "E12" '  +2.84%  '

Set-TextValue "D13" '19.98'
Set-TextValue "E13" '  +0.17%  '

Set-TextValue "E14" '  +0.16%  '

Set-TextValue "D15" '3.348.81'
Set-TextValue "E15" '  +3.66%  '

Set-TextValue "D16" '0.997'
Set-TextValue "E16" '  +6.83%  '

Set-TextValue "D17" '2.908.87'
Set-TextValue "E17" '  +3.61%  '

Set-TextValue "D18" '52.140.86'
Set-TextValue "E18" '  +1.14%  '

Set-TextValue "D19" '7.72'
Set-TextValue "E19" '  +0.18%  '

Set-TextValue "E20" '  +6.46%  '

Set-TextValue "D21" '14.36'
Set-TextValue "E21" '  +7.84%  '

Set-TextValue "D22" '0.0₃0980'
Set-TextValue "E22" '  +1.18%  '

Set-TextValue "D23" '71.01'
Set-TextValue "E23" '  +0.97%  '

Set-TextValue "D24" '270.24'
Set-TextValue "E24" '  +1.35%  '

Set-TextValue "E25" '  +1.48%  '

Set-TextValue "D26" '26.48'
Set-TextValue "E26" '  +2.18%  '

Set-TextValue "E28" '  -0.67%  '

Set-TextValue "D29" '38.75'
Set-TextValue "E29" '  +4.58%  '

Set-TextValue "D30" '10.52'
Set-TextValue "E30" '  +2.15%  '

Set-TextValue "E31" '  +1.09%  '

Set-TextValue "D32" '6.47'
Set-TextValue "E32" '  +3.47%  '

Set-TextValue "E33" '  +9.10%  '

Set-TextValue "B34" 'Hedera'
Set-TextValue "C34" 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue "D34" '0.0949'
Set-TextValue "E34" '  +11.69%  '

Set-TextValue "B35" 'OKB'
Set-TextValue "C35" 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue "D35" '53.26'
Set-TextValue "E35" '  +1.97%  '

Set-TextValue "E36" '  +4.10%  '

Set-TextValue "D37" '0.998'
Set-TextValue "E37" '  -0.18%  '

Set-TextValue "E38" '  +6.08%  '

Set-TextValue "D39" '18.73'
Set-TextValue "E39" '  +1.01%  '

Set-TextValue "E40" '  +3.78%  '

Set-TextValue "E41" '  +5.98%  '

Set-TextValue "E42" '  +2.53%  '

Set-TextValue "D43" '22.77'
Set-TextValue "E43" '  +3.85%  '

Set-TextValue "E44" '  +1.59%  '

Set-TextValue "E45" '  +0.85%  '

Set-TextValue "E46" '  +6.87%  '

Set-TextValue "D47" '2.203.99'
Set-TextValue "E47" '  +3.45%  '

Set-TextValue "D48" '2.51'
Set-TextValue "E48" '  +6.82%  '

Set-TextValue "E49" '  +19.16%  '

Set-TextValue "D50" '0.950'
Set-TextValue "E50" '  +4.17%  '

Set-TextValue "D51" '5.52'
Set-TextValue "E51" '  +3.43%  '
